$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp update ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 30 de Mayo de 2020 a las 06:40"

# --- Re-sorted / updated country rows ---
# Kirguistan moves above Mayotte; Kirguistan gets fresh case counts,
# Mayotte keeps its previous numbers (shifted down one row).
$ws.Cells.Item(96,1).Value = "Kirguistan"
$ws.Cells.Item(96,2).Value = 1722
$ws.Cells.Item(96,3).Value = 60
$ws.Cells.Item(96,4).Value = 1113
$ws.Cells.Item(96,5).Value = 593
$ws.Cells.Item(96,6).Value = 0
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 16

$ws.Cells.Item(97,1).Value = "Mayotte"
$ws.Cells.Item(97,2).Value = 1699
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 1385
$ws.Cells.Item(97,5).Value = 293
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 21

# Curazao moves above Fiyi; Santa Lucia moves above Belice.
$ws.Cells.Item(198,1).Value = "Curazao"
$ws.Cells.Item(198,2).Value = 18
$ws.Cells.Item(198,3).Value = 0
$ws.Cells.Item(198,4).Value = 14
$ws.Cells.Item(198,5).Value = 3
$ws.Cells.Item(198,6).Value = 0
$ws.Cells.Item(198,7).Value = 0
$ws.Cells.Item(198,8).Value = 1

$ws.Cells.Item(199,1).Value = "Fiyi"
$ws.Cells.Item(199,2).Value = 18
$ws.Cells.Item(199,3).Value = 0
$ws.Cells.Item(199,4).Value = 15
$ws.Cells.Item(199,5).Value = 3
$ws.Cells.Item(199,6).Value = 0
$ws.Cells.Item(199,7).Value = 0
$ws.Cells.Item(199,8).Value = 0

$ws.Cells.Item(200,1).Value = "Santa Lucia"
$ws.Cells.Item(200,2).Value = 18
$ws.Cells.Item(200,3).Value = 0
$ws.Cells.Item(200,4).Value = 18
$ws.Cells.Item(200,5).Value = 0
$ws.Cells.Item(200,6).Value = 0
$ws.Cells.Item(200,7).Value = 0
$ws.Cells.Item(200,8).Value = 0

$ws.Cells.Item(201,1).Value = "Belice"
$ws.Cells.Item(201,2).Value = 18
$ws.Cells.Item(201,3).Value = 0
$ws.Cells.Item(201,4).Value = 16
$ws.Cells.Item(201,5).Value = 0
$ws.Cells.Item(201,6).Value = 0
$ws.Cells.Item(201,7).Value = 0
$ws.Cells.Item(201,8).Value = 2

# San Bartolome moves above Bonaire, San Eustaquio y Saba (values unchanged).
$ws.Cells.Item(215,1).Value = "San Bartolome"
$ws.Cells.Item(215,2).Value = 6
$ws.Cells.Item(215,3).Value = 0
$ws.Cells.Item(215,4).Value = 6
$ws.Cells.Item(215,5).Value = 0
$ws.Cells.Item(215,6).Value = 0
$ws.Cells.Item(215,7).Value = 0
$ws.Cells.Item(215,8).Value = 0

$ws.Cells.Item(216,1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216,2).Value = 6
$ws.Cells.Item(216,3).Value = 0
$ws.Cells.Item(216,4).Value = 6
$ws.Cells.Item(216,5).Value = 0
$ws.Cells.Item(216,6).Value = 0
$ws.Cells.Item(216,7).Value = 0
$ws.Cells.Item(216,8).Value = 0
